$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Changelog block added below the existing notes (rows 21-23) ---------
$ws.Range("A21").Value = "Changelog:"
$ws.Range("A22").Value = "Juli 26"

# The EU-cord note is entered before the US-cord note so the shared string
# table ends up with the same ordering as the saved workbook (EU text gets
# the lower shared-string index even though it is written into the later
# row, C23, while the US text lands in the earlier row, C22).
$ws.Range("C23").Value = "The partnumber for the EU cord has been added"
$ws.Range("C22").Value = "The partnumber for the US cord has been corrected"

# Give the two new note cells the same left/centre alignment used by the
# existing note in C19 (vertical alignment is applied first so the engine
# settles directly on the combined horizontal+vertical style).
$ws.Range("C22").VerticalAlignment = -4108
$ws.Range("C22").HorizontalAlignment = -4131
$ws.Range("C23").VerticalAlignment = -4108
$ws.Range("C23").HorizontalAlignment = -4131

# Matches the saved workbook's final selection.
[void]$ws.Range("H21").Select()
